$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(585, 331)'
$ws.Range("B2").Value = '(450, 321)'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '135.36986370680884'
$ws.Range("D2").Value = '0:00:00.591743'
$ws.Range("E2").Value = '[[553, 327, datetime.timedelta(microseconds=94817), 653.527785899376], [516, 325, datetime.timedelta(microseconds=202944), 84.83404513007834, -2802.22002507735], [481, 317, datetime.timedelta(microseconds=317487), 52.54260911975261, -101.7094747511732], [460, 319, datetime.timedelta(microseconds=413743), 23.689749565895667, -69.73618781189518], [450, 323, datetime.timedelta(microseconds=503781), 9.933409370623474, -27.306190974396007]]'
$ws.Range("F2").Value = '2022-07-18 16:31:25.567586'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '62.89746106377338'
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = '51.73973328110612'
$ws.Range("A3").Value = '(672, 350)'
$ws.Range("B3").Value = '(823, 389)'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '155.95512174981621'
$ws.Range("D3").Value = '0:00:00.620233'
$ws.Range("E3").Value = '[[700, 356, datetime.timedelta(microseconds=96082), 79.75419590077422], [739, 364, datetime.timedelta(microseconds=188565), 98.0990205156307, 97.28647742081765], [764, 364, datetime.timedelta(microseconds=305908), 37.97172240849074, -196.55353278482406], [799, 379, datetime.timedelta(microseconds=418501), 42.27646093764834, 10.286088991800728], [823, 391, datetime.timedelta(microseconds=511787), 24.36060789087095, -35.00646371786972]]'
$ws.Range("F3").Value = '2022-07-18 16:31:31.739858'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '72.46207486180485'
$ws.Range("I3").Value = '(864.0, 360.0)'
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = '50.21951811795888'
$ws.Range("A4").Value = '(723, 350)'
$ws.Range("B4").Value = '(827, 385)'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '109.73149046650191'
$ws.Range("D4").Value = '0:00:00.484416'
$ws.Range("E4").Value = '[[778, 362, datetime.timedelta(microseconds=95326), 460.75635450464534], [825, 364, datetime.timedelta(microseconds=196310), 111.34209974754187, -1779.9106248133235], [823, 393, datetime.timedelta(microseconds=291252), 46.373573285444444, -223.0663702295518], [827, 385, datetime.timedelta(microseconds=386315), 10.75757903148568, -92.19417898336529]]'
$ws.Range("F4").Value = '2022-07-18 16:31:38.012506'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '50.98499739967954'
$ws.Range("I4").Value = '(864.0, 360.0)'
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = '44.654227123532216'
$ws.Range("A5").Value = '(684, 352)'
$ws.Range("B5").Value = '(844, 397)'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '166.20770138594662'
$ws.Range("D5").Value = '0:00:00.606460'
$ws.Range("E5").Value = '[[715, 366, datetime.timedelta(microseconds=113997), 138.6386690320643], [751, 352, datetime.timedelta(microseconds=207464), 86.50730606838702, -251.27907956887603], [784, 395, datetime.timedelta(microseconds=304395), 82.73694952238448, -12.386394474293422], [817, 397, datetime.timedelta(microseconds=403527), 38.06699592250878, -110.69879735402017], [842, 397, datetime.timedelta(microseconds=505468), 22.980393731228453, -29.846799780164773]]'
$ws.Range("F5").Value = '2022-07-18 16:31:44.202323'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '77.22577344883618'
$ws.Range("I5").Value = '(864.0, 360.0)'
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = '42.05948168962618'

Write-Output "done"